$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Neo4j query text for the "SamplesTab" row (B3): now reads the tumor
# status directly from the sample node instead of collecting it separately.
$newSamplesQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Molecular Characterization Initiative"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newSamplesQuery

# The row grew taller to fit the updated (slightly longer) query text.
$ws.Rows(3).RowHeight = 218.25

# Move / record the active selection as it was left after the edit.
$ws.Range("B11").Select()
